$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.Value = "'67.675.30"
$c.ClearFormats()
$ws.Range("E2").Value = "  +2.67%  "

$c = $ws.Range("D3")
$c.Value = "'3.317.08"
$c.ClearFormats()
$ws.Range("E3").Value = "  +0.07%  "

$ws.Range("E4").Value = "  -0.16%  "

$c = $ws.Range("D5")
$c.Value = "'584.89"
$c.ClearFormats()
$ws.Range("E5").Value = "  +4.52%  "

$c = $ws.Range("D6")
$c.Value = "'182.09"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.50%  "

$c = $ws.Range("D7")
$c.Value = "'1.00"
$c.ClearFormats()
$ws.Range("E7").Value = "  +0.22%  "

$c = $ws.Range("D8")
$c.Value = "'0.590"
$c.ClearFormats()
$ws.Range("E8").Value = "  +3.32%  "

$c = $ws.Range("D9")
$c.Value = "'3.313.09"
$c.ClearFormats()
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("E10").Value = "  +1.25%  "

$c = $ws.Range("D11")
$c.Value = "'0.578"
$c.ClearFormats()
$ws.Range("E11").Value = "  +0.44%  "

$c = $ws.Range("D12")
$c.Value = "'46.21"
$c.ClearFormats()
$ws.Range("E12").Value = "  +1.07%  "

$ws.Range("E13").Value = "  +3.71%  "

$c = $ws.Range("D14")
$c.Value = "'641.21"
$c.ClearFormats()
$ws.Range("E14").Value = "  +10.83%  "

$c = $ws.Range("D15")
$c.Value = "'3.851.70"
$c.ClearFormats()
$ws.Range("E15").Value = "  +0.27%  "

$c = $ws.Range("D16")
$c.Value = "'8.43"
$c.ClearFormats()
$ws.Range("E16").Value = "  +0.32%  "

$c = $ws.Range("D17")
$c.Value = "'67.766.04"
$c.ClearFormats()
$ws.Range("E17").Value = "  +2.80%  "

$c = $ws.Range("D19")
$c.Value = "'3.315.43"
$c.ClearFormats()
$ws.Range("E19").Value = "  +0.37%  "

$c = $ws.Range("D20")
$c.Value = "'17.62"
$c.ClearFormats()
$ws.Range("E20").Value = "  +0.36%  "

$c = $ws.Range("D21")
$c.Value = "'10.88"
$c.ClearFormats()
$ws.Range("E21").Value = "  +0.12%  "

$c = $ws.Range("D22")
$c.Value = "'0.900"
$c.ClearFormats()
$ws.Range("E22").Value = "  +1.08%  "

$c = $ws.Range("D23")
$c.Value = "'17.61"
$c.ClearFormats()
$ws.Range("E23").Value = "  -0.01%  "

$c = $ws.Range("D24")
$c.Value = "'5.03"
$c.ClearFormats()
$ws.Range("E24").Value = "  +0.87%  "

$c = $ws.Range("D25")
$c.Value = "'97.08"
$c.ClearFormats()
$ws.Range("E25").Value = "  -0.76%  "

$c = $ws.Range("D26")
$c.Value = "'4.01"
$c.ClearFormats()
$ws.Range("E26").Value = "  +1.77%  "

$ws.Range("E27").Value = "  +3.15%  "

$c = $ws.Range("D28")
$c.Value = "'9.55"
$c.ClearFormats()
$ws.Range("E28").Value = "  +2.59%  "

$ws.Range("E29").Value = "  +7.09%  "

$c = $ws.Range("D30")
$c.Value = "'8.55"
$c.ClearFormats()
$ws.Range("E30").Value = "  +1.32%  "

$c = $ws.Range("D31")
$c.Value = "'6.66"
$c.ClearFormats()
$ws.Range("E31").Value = "  +0.65%  "

$c = $ws.Range("D32")
$c.Value = "'589.72"
$c.ClearFormats()
$ws.Range("E32").Value = "  +5.51%  "

$c = $ws.Range("D33")
$c.Value = "'3.930.37"
$c.ClearFormats()
$ws.Range("E33").Value = "  +5.11%  "

$c = $ws.Range("D34")
$c.Value = "'10.93"
$c.ClearFormats()
$ws.Range("E34").Value = "  +1.52%  "

$ws.Range("E35").Value = "  -4.50%  "

$ws.Range("E36").Value = "  +1.39%  "

$ws.Range("E37").Value = "  -0.25%  "

$c = $ws.Range("D38")
$c.Value = "'55.51"
$c.ClearFormats()
$ws.Range("E38").Value = "  -0.22%  "

$c = $ws.Range("D39")
$c.Value = "'0.129"
$c.ClearFormats()
$ws.Range("E39").Value = "  +2.14%  "

$ws.Range("E40").Value = "  +3.47%  "

$ws.Range("E41").Value = "  +3.95%  "

$c = $ws.Range("D42")
$c.Value = "'32.52"
$c.ClearFormats()
$ws.Range("E42").Value = "  -2.08%  "

$c = $ws.Range("D43")
$c.Value = "'3.39"
$c.ClearFormats()
$ws.Range("E43").Value = "  +1.77%  "

$c = $ws.Range("D44")
$c.Value = "'0.0₃0683"
$c.ClearFormats()
$ws.Range("E44").Value = "  +0.34%  "

$c = $ws.Range("D45")
$c.Value = "'0.337"
$c.ClearFormats()
$ws.Range("E45").Value = "  +1.74%  "

$c = $ws.Range("D46")
$c.Value = "'0.0413"
$c.ClearFormats()
$ws.Range("E46").Value = "  +1.60%  "

$ws.Range("E47").Value = "  +1.55%  "

$c = $ws.Range("D48")
$c.Value = "'1.00"
$c.ClearFormats()
$ws.Range("E48").Value = "  +0.56%  "

$ws.Range("B49").Value = "ThetaToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$c = $ws.Range("D49")
$c.Value = "'2.54"
$c.ClearFormats()
$ws.Range("E49").Value = "  +1.91%  "

$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D50")
$c.Value = "'1.37"
$c.ClearFormats()
$ws.Range("E50").Value = "  +12.44%  "

$c = $ws.Range("D51")
$c.Value = "'130.53"
$c.ClearFormats()
$ws.Range("E51").Value = "  +4.14%  "
